$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 15 and 16 swap (Chainlink <-> WrappedEther) plus value updates
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.46"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.20%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.514.98"
$ws.Range("E16").Value = "  -0.23%  "

# Update Price (D) and Volume(1h) (E) columns for remaining rows
$ws.Range("D2").Value = "42.650.18"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").Value = "2.530.60"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.39"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.23"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -5.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.575"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.05%  "
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("E9").Value = "  -3.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.19"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0808"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.55"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.88%  "
$ws.Range("E13").Value = "  -0.44%  "
$ws.Range("D14").Value = "2.919.15"
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("D18").Value = "42.624.88"
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.98"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.33%  "
$ws.Range("E20").Value = "  +0.96%  "
$ws.Range("E21").Value = "  -2.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.05"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "250.88"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.29%  "
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.00"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.61"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.38%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.40"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.28"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.67%  "
$ws.Range("E30").Value = "  -3.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.04"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "154.72"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.35%  "
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.13"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.61%  "
$ws.Range("E35").Value = "  -1.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0785"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.82%  "
$ws.Range("E37").Value = "  -0.42%  "
$ws.Range("E38").Value = "  -4.33%  "
$ws.Range("E39").Value = "  -1.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.69"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.31"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +10.75%  "
$ws.Range("E42").Value = "  -3.11%  "
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("E44").Value = "  -1.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.29"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -5.22%  "
$ws.Range("D46").Value = "2.016.27"
$ws.Range("E46").Value = "  -1.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "85.50"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.80"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.20%  "
$ws.Range("D49").Value = "2.773.31"
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.94"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.80%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "102.50"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.16%  "
